# Updated R3 and Time_Worked
# This script applies the "Andreas removed / Research section added / hours updated"
# edit described by the commit to Reports/Time_Worked.xlsx (sheet1 = Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Remove the "Andreas" row from every existing table. ClearContents()
#    drops the value but keeps whatever direct formatting (s=...) the cell
#    already carried, and cells that had no explicit style disappear
#    entirely - matching the target XML for rows 7, 15, 23, 31, 39.
#    Doing this first also lets the now-unused "Andreas" shared string get
#    garbage collected before we add the new "Research" string, so
#    "Research" lands in the freed slot (matches uniqueCount staying 15).
# ---------------------------------------------------------------------------
$ws.Range("A7:F7").ClearContents()
$ws.Range("A15:F15").ClearContents()
$ws.Range("A23:F23").ClearContents()
$ws.Range("A31:F31").ClearContents()
$ws.Range("A39:F39").ClearContents()

# ---------------------------------------------------------------------------
# 2) Update the logged hours in the existing tables (Wednesday/Thursday
#    columns mostly gain real hours instead of 0, picking up the
#    highlighted "s=2"/"s=3" look already used elsewhere in the sheet).
# ---------------------------------------------------------------------------

# -- "Total" table (rows 1-7): D3:E6 gain hours and the highlighted style --
$ws.Range("B3").Copy()
$ws.Range("D3:E6").PasteSpecial(-4122)
$ws.Range("D3").Value = 4.1
$ws.Range("E3").Value = 5.1
$ws.Range("D4").Value = 3.6
$ws.Range("E4").Value = 4.6
$ws.Range("D5").Value = 4.1
$ws.Range("E5").Value = 6.1
$ws.Range("D6").Value = 3.1
$ws.Range("E6").Value = 5.1

# -- "Documentation" table (rows 9-15) --
$ws.Range("B3").Copy()
$ws.Range("D11:E11").PasteSpecial(-4122)
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("D13:E14").PasteSpecial(-4122)
$ws.Range("D11").Value = 0.5
$ws.Range("E11").Value = 1.5
$ws.Range("E12").Value = 1
$ws.Range("D13").Value = 0.5
$ws.Range("E13").Value = 1.5
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 1.5

# -- "UML" table (rows 33-39): D35:D38 take the fill-only style (s=3, like
#    B35/B37/B38), E35:E38 take the highlighted style (s=2) --
$ws.Range("B35").Copy()
$ws.Range("D35:D38").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("E35:E38").PasteSpecial(-4122)
$ws.Range("D35").Value = 3.6
$ws.Range("E35").Value = 3.6
$ws.Range("D36").Value = 3.6
$ws.Range("E36").Value = 3.6
$ws.Range("D37").Value = 3.6
$ws.Range("E37").Value = 3.6
$ws.Range("D38").Value = 2.6
$ws.Range("E38").Value = 3.6

# ---------------------------------------------------------------------------
# 3) Add the new "Research" table (rows 41-46), mirroring the structure of
#    the other small tables (label row, header row, 4 employee rows - no
#    "Andreas" row since that name is gone from the workbook).
# ---------------------------------------------------------------------------
$ws.Range("C41").Value = "Research"

$ws.Range("B42").Value = "Monday"
$ws.Range("C42").Value = "Tuesday"
$ws.Range("D42").Value = "Wednesday"
$ws.Range("E42").Value = "Thursday"
$ws.Range("F42").Value = "Friday"

$ws.Range("A43").Value = "Rasmus"
$ws.Range("A44").Value = "Erik"
$ws.Range("A45").Value = "Kim"
$ws.Range("A46").Value = "Calle"

# Give B43:F46 the same plain right-aligned style (s=1) used by the other
# tables' data rows ...
$ws.Range("B19:F22").Copy()
$ws.Range("B43:F46").PasteSpecial(-4122)

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0

$ws.Range("B44").Value = 0
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0

$ws.Range("B45").Value = 0
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 1
$ws.Range("F45").Value = 0

$ws.Range("B46").Value = 0
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0

# ... except E45 (Kim/Thursday = 1 hour), which gets the highlighted style.
$ws.Range("B3").Copy()
$ws.Range("E45").PasteSpecial(-4122)
$ws.Range("E45").Value = 1

# ---------------------------------------------------------------------------
# 4) Restore the view: scroll so row 16 is at the top, and leave the
#    selection where the author left it (I45).
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("I45").Select()
$excel.ActiveWindow.ScrollRow = 16

$wb.Application.CutCopyMode = $false
